$d = $word.ActiveDocument

# The template has several "{%- ... %}" tags; the only one whose full
# paragraph text is exactly "{%- endfor %}" is the `for ach in
# job.achievements` loop terminator right after the "{{ ach }}" line.
# That text string is unique in the whole document, so Find.Execute on
# it unambiguously locates the right paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("{%- endfor %}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate '{%- endfor %}' in the document"
}

# $rng now spans the matched text "{%- endfor %}". The hyphen that
# forms the Jinja "trim" marker sits right after the opening "{%",
# i.e. at offset 2 from the start of the match. Remove just that one
# character so the run holding "{%" and the run holding the trailing
# " " collapse back together into a single "{% " run, exactly like the
# rest of the already-merged "{% " tags elsewhere in the document -
# leaving the "endfor" run/proofErr markers and the " %}" run
# completely untouched.
$dashRange = $d.Range($rng.Start + 2, $rng.Start + 3)
if ($dashRange.Text -ne "-") {
    throw "Unexpected character where the trim-hyphen was expected: [$($dashRange.Text)]"
}
$dashRange.Delete()
